# RP3 Flights 2021 Jan-Dec — "2021 post ops updates"
#
# Updates the ERT_FLTS_YY sheet (the single source-of-truth sheet; the
# other sheets pull the release date via formula and already hold the
# correct cumulative totals) with:
#   - a new release date (B2)
#   - full-year day counts for every year row (B6:B12), replacing the
#     prior partial-year day counts. D (Avg. Daily) and E (% change) are
#     formulas and recompute automatically, as do the cross-sheet B2
#     formulas on ERT_FLTS_MM / ERT_FLTS_LOC.
#
# Column widths also shrank uniformly (~x0.875) across every sheet in the
# source edit (a side effect of the resave, not a deliberate per-cell
# resize); ColumnWidth is reproduced here too, snapped to the nearest
# value this engine's pixel-grid rounding can represent.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# ERT_FLTS_YY
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ERT_FLTS_YY")

$ws.Range("B2").Value = 44665

$ws.Range("B6").Value  = 365   # 2015
$ws.Range("B7").Value  = 366   # 2016 (leap)
$ws.Range("B8").Value  = 365   # 2017
$ws.Range("B9").Value  = 365   # 2018
$ws.Range("B10").Value = 365   # 2019
$ws.Range("B11").Value = 366   # 2020 (leap)
$ws.Range("B12").Value = 365   # 2021

$ws.Columns.Item(1).ColumnWidth = 10.666666666666666
$ws.Columns.Item(2).ColumnWidth = 16.0
$ws.Columns.Item(3).ColumnWidth = 8.333333333333334
$ws.Columns.Item(4).ColumnWidth = 9.0
$ws.Columns.Item(5).ColumnWidth = 8.0
$ws.Columns.Item(6).ColumnWidth = 10.666666666666666

# ---------------------------------------------------------------------
# ERT_FLTS_MM
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ERT_FLTS_MM")

$ws.Columns.Item(1).ColumnWidth = 10.666666666666666
$ws.Columns.Item(2).ColumnWidth = 15.833333333333334
$ws.Columns.Item(3).ColumnWidth = 14.333333333333334
$ws.Columns.Item(4).ColumnWidth = 11.333333333333334
$ws.Columns.Item(5).ColumnWidth = 8.333333333333334
$ws.Columns.Item(6).ColumnWidth = 8.333333333333334
$ws.Columns.Item(7).ColumnWidth = 7.333333333333333
$ws.Columns.Item(8).ColumnWidth = 9.666666666666666
$ws.Columns.Item(9).ColumnWidth = 6.833333333333333

# ---------------------------------------------------------------------
# ERT_FLTS_LOC
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ERT_FLTS_LOC")

$ws.Columns.Item(1).ColumnWidth = 16.833333333333332
$ws.Columns.Item(2).ColumnWidth = 16.0
$ws.Columns.Item(3).ColumnWidth = 11.0
$ws.Columns.Item(4).ColumnWidth = 12.0
$ws.Columns.Item(5).ColumnWidth = 12.5
$ws.Columns.Item(6).ColumnWidth = 18.333333333333332

# ---------------------------------------------------------------------
# Change Log
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Change Log")

$ws.Columns.Item(1).ColumnWidth = 11.333333333333334
$ws.Columns.Item(2).ColumnWidth = 8.666666666666666
$ws.Columns.Item(3).ColumnWidth = 11.833333333333334
$ws.Columns.Item(4).ColumnWidth = 122.33333333333333
